$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Create WO")

# Update B2 value: "Pro-Lot Track (Lot Track)" -> "Pro-SYDATA1 (Lot track)"
$ws.Range("B2").Value = "Pro-SYDATA1 (Lot track)"

# Remove the bold header styling (font) applied to A1:F1, resetting to default style
$ws.Range("A1:F1").ClearFormats()

# Update the selection to match A1:K2 (active cell lands on B2 per the target view state)
$ws.Range("A1:K2").Select()
